$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 conversion text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.76 = 40965.07 pesos`n✅ 40965.07 pesos = 9.68 = 942.83 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $text

# --- Update tasas sheet N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 102.48
$ws2.Range("O10").Value = 4198.1
$ws2.Range("N12").Value = 4229.88
$ws2.Range("O12").Value = 97.35299999999999
